$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7763
$ws.Range("F3").Value = 3643
$ws.Range("F9").Value = 195
$ws.Range("F11").Value = 35
$ws.Range("F14").Value = 4386
$ws.Range("F15").Value = 4386
$ws.Range("F17").Value = 435
$ws.Range("F18").Value = 1056
$ws.Range("F20").Value = 3794
$ws.Range("F21").Value = 136
$ws.Range("F22").Value = 127
$ws.Range("F23").Value = 24
$ws.Range("F24").Value = 141
$ws.Range("F25").Value = 3208
$ws.Range("F26").Value = 2602
$ws.Range("F30").Value = 113
$ws.Range("F33").Value = 58
$ws.Range("F36").Value = 46
$ws.Range("F38").Value = 4823
$ws.Range("F39").Value = 615
$ws.Range("F43").Value = 937
$ws.Range("F44").Value = 334
$ws.Range("F45").Value = 17
$ws.Range("F46").Value = 1796
$ws.Range("F49").Value = 651
$ws.Range("F50").Value = 775

$ws = $wb.Worksheets.Item(2)
$ws.Range("F23").Value = 683

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 7763
$ws.Range("F5").Value = 3643
$ws.Range("F11").Value = 195
$ws.Range("F13").Value = 35
$ws.Range("F15").Value = 4386
$ws.Range("F16").Value = 4386
$ws.Range("F20").Value = 435
$ws.Range("F21").Value = 1056
$ws.Range("F23").Value = 3794
$ws.Range("F24").Value = 136
$ws.Range("F25").Value = 127
$ws.Range("F26").Value = 3208
$ws.Range("F27").Value = 2602
$ws.Range("F30").Value = 113
$ws.Range("F33").Value = 58
$ws.Range("F36").Value = 46
$ws.Range("F39").Value = 4824
$ws.Range("F41").Value = 615
$ws.Range("F45").Value = 937
$ws.Range("F46").Value = 334
$ws.Range("F47").Value = 1796
$ws.Range("F49").Value = 651
$ws.Range("F50").Value = 775
